# Updated cryptos list data (price / 1h volume figures, and a row-order
# correction for ShibaInu / WrappedEther) as per the latest GitHub Actions run.
#
# Note: several "Price" values look like plain numbers (e.g. "523.67").
# Excel's COM layer auto-converts such strings to numeric cell values,
# which would silently drop formatting like trailing zeros / leading
# apostrophe-worthy ambiguity (e.g. "0.420" -> 0.42, "1.00" -> 1,
# "0.0000136" -> 1.36E-05). To keep these as text (matching the original
# inline-string storage), a leading apostrophe is prepended for those
# values before assigning - exactly what typing '523.67 into a cell does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16/17: ShibaInu and WrappedEther swapped positions, each with fresh
# price / volume data.
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.664.85"
$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = '''0.0000136'
$ws.Range("E17").Value = "  -2.33%  "

# Row 2: Bitcoin
$ws.Range("D2").Value = '58.920.51'
$ws.Range("E2").Value = '  -3.27%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.660.41'
$ws.Range("E3").Value = '  -1.47%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.04%  '

# Row 5: BNB
$ws.Range("D5").Value = '''523.67'
$ws.Range("E5").Value = '  -0.14%  '

# Row 6: Solana
$ws.Range("D6").Value = '''144.19'
$ws.Range("E6").Value = '  -2.22%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.07%  '

# Row 8: XRP
$ws.Range("D8").Value = '''0.569'
$ws.Range("E8").Value = '  -1.54%  '

# Row 9: Toncoin
$ws.Range("E9").Value = '  +7.45%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.102'
$ws.Range("E10").Value = '  -3.72%  '

# Row 11: Cardano
$ws.Range("E11").Value = '  -2.01%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +1.33%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.128.49'
$ws.Range("E13").Value = '  -2.11%  '

# Row 14: WrappedBTC
$ws.Range("D14").Value = '58.920.16'
$ws.Range("E14").Value = '  -2.96%  '

# Row 15: Avalanche
$ws.Range("D15").Value = '''20.99'
$ws.Range("E15").Value = '  -2.10%  '

# Row 18: BitcoinCash
$ws.Range("D18").Value = '''338.55'
$ws.Range("E18").Value = '  -4.23%  '

# Row 19: Polkadot
$ws.Range("E19").Value = '  -3.70%  '

# Row 20: Chainlink
$ws.Range("E20").Value = '  -2.62%  '

# Row 21: Uniswap
$ws.Range("E21").Value = '  -0.04%  '

# Row 22: Dai
$ws.Range("E22").Value = '  -0.07%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''64.56'
$ws.Range("E23").Value = '  +2.43%  '

# Row 24: Polygon
$ws.Range("D24").Value = '''0.420'

# Row 25: Kaspa
$ws.Range("E25").Value = '  -1.68%  '

# Row 26: Binance-PegBSC-USD
$ws.Range("D26").Value = '''0.997'
$ws.Range("E26").Value = '  -0.46%  '

# Row 27: PEPE
$ws.Range("E27").Value = '  -2.87%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = '''7.14'
$ws.Range("E28").Value = '  -2.61%  '

# Row 29: Aptos
$ws.Range("D29").Value = '''6.66'
$ws.Range("E29").Value = '  -3.22%  '

# Row 30: USDe
$ws.Range("E30").Value = '  -0.04%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -0.56%  '

# Row 32: EthereumClassic
$ws.Range("E32").Value = '  -1.42%  '

# Row 33: Monero
$ws.Range("D33").Value = '''150.41'
$ws.Range("E33").Value = '  +1.81%  '

# Row 34: NEARProtocol
$ws.Range("E34").Value = '  -4.30%  '

# Row 35: ImmutableX
$ws.Range("E35").Value = '  -5.45%  '

# Row 36: SuiNetwork
$ws.Range("E36").Value = '  -6.14%  '

# Row 37: Fetch.AI
$ws.Range("E37").Value = '  -1.60%  '

# Row 38: OKB
$ws.Range("D38").Value = '''36.82'
$ws.Range("E38").Value = '  -0.44%  '

# Row 39: Stacks
$ws.Range("E39").Value = '  -5.81%  '

# Row 40: Filecoin
$ws.Range("E40").Value = '  -3.48%  '

# Row 41: Mantle
$ws.Range("D41").Value = '''0.615'
$ws.Range("E41").Value = '  -0.01%  '

# Row 42: FirstDigitalUSD
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.25%  '

# Row 43: Bittensor
$ws.Range("D43").Value = '''275.11'
$ws.Range("E43").Value = '  -3.86%  '

# Row 44: EnergySwap
$ws.Range("D44").Value = '''19.81'
$ws.Range("E44").Value = '  -1.75%  '

# Row 45: Stellar
$ws.Range("E45").Value = '  -2.34%  '

# Row 46: WhiteBITCoin
$ws.Range("D46").Value = '''10.67'
$ws.Range("E46").Value = '  +2.02%  '

# Row 47: Hedera
$ws.Range("E47").Value = '  -1.51%  '

# Row 48: Maker
$ws.Range("D48").Value = '2.051.63'
$ws.Range("E48").Value = '  -4.53%  '

# Row 49: RenderToken
$ws.Range("D49").Value = '''4.71'
$ws.Range("E49").Value = '  -3.39%  '

# Row 50: VeChain
$ws.Range("E50").Value = '  -3.13%  '

# Row 51: InjectiveProtocol
$ws.Range("D51").Value = '''18.81'
$ws.Range("E51").Value = '  -3.53%  '
